$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# SAMPLEONE (sheet2): replace the generic placeholder Name/Item labels with
# the real character / organisation names, and bump a few sample numbers.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("SAMPLEONE")

$ws1.Range("D5").Value = "御坂美琴"
$ws1.Range("G5").Value = "LEVEL5"
$ws1.Range("K5").Value = "常盘台中学"
$ws1.Range("O5").Value = "第三位"

$ws1.Range("D6").Value = "白井黑子"
$ws1.Range("G6").Value = "LEVEL4"
$ws1.Range("H6").Value = 27
$ws1.Range("K6").Value = "177支部所属风纪委员"
$ws1.Range("L6").Value = 34

$ws1.Range("D7").Value = "初春饰利"
$ws1.Range("G7").Value = "LEVEL1"
$ws1.Range("H7").Value = 29

$ws1.Range("D8").Value = "佐天泪子"

$ws1.Range("D9").Value = "食蜂 操祈"

# ---------------------------------------------------------------------------
# SAMPLEFOUR (sheet5): rename the "float64" type label to "float".
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("SAMPLEFOUR")

$ws4.Range("I2").Value = "float"
$ws4.Range("M2").Value = "float"
$ws4.Range("Q2").Value = "float"

# ---------------------------------------------------------------------------
# Selection / active-tab bookkeeping: SAMPLEFOUR loses the tab focus, it
# moves to SAMPLEONE, and each sheet keeps its own remembered selection.
# ---------------------------------------------------------------------------
[void]$ws4.Activate()
[void]$ws4.Range("Q2").Select()

[void]$ws1.Activate()
[void]$ws1.Range("G10").Select()
